# Q23018277-en.xlsx: the two "secondary source" rows (Joe Biden article and
# Office of the Clerk article) swapped places in the underlying data pipeline
# (e.g. the directory listing order used to build this sheet shifted after a
# new JSON file was added elsewhere in the dataset), so row 3 now holds what
# used to be row 4's data and vice versa. Columns C ("historical distance")
# and D ("time bucket") are "unknown" in both rows, so they are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for row 3 and row 4.
$titleRow3 = $ws.Range("A3").Value()
$timestampRow3 = $ws.Range("B3").Value()
$uriRow3 = $ws.Range("E3").Value()

$titleRow4 = $ws.Range("A4").Value()
$timestampRow4 = $ws.Range("B4").Value()
$uriRow4 = $ws.Range("E4").Value()

# Write row 4's former content into row 3 ...
$ws.Range("A3").Value = $titleRow4
$ws.Range("B3").Value = $timestampRow4
$ws.Range("E3").Value = $uriRow4

# ... and row 3's former content into row 4.
$ws.Range("A4").Value = $titleRow3
$ws.Range("B4").Value = $timestampRow3
$ws.Range("E4").Value = $uriRow3
